# Fruta / hortaliza, semanal
# Update the weekly Jengibre price records (rows 4-13) so that each row
# now holds the data previously found in the following row, and the
# last row receives a brand new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 4;  D = 44764; J = 200; K = 12000; L = 13000; M = 12500; P = 962 },
    @{ Row = 5;  D = 44389; J = 120; K = 12000; L = 13000; M = 12500; P = 962 },
    @{ Row = 6;  D = 44320; J = 160; K = 19000; L = 20000; M = 19500; P = 1500 },
    @{ Row = 7;  D = 44379; J = 120; K = 12000; L = 13000; M = 12667; P = 974 },
    @{ Row = 8;  D = 44580; J = 160; K = 11000; L = 12000; M = 11500; P = 885 },
    @{ Row = 9;  D = 44397; J = 140; K = 12500; L = 13000; M = 12750; P = 981 },
    @{ Row = 10; D = 44592; J = 120; K = 12000; L = 13000; M = 12500; P = 962 },
    @{ Row = 11; D = 44159; J = 100; K = 23000; L = 24000; M = 23500; P = 1808 },
    @{ Row = 12; D = 44406; J = 160; K = 17000; L = 18000; M = 17500; P = 1346 },
    @{ Row = 13; D = 44832; J = 100; K = 13000; L = 14000; M = 13500; P = 1038 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("D$r").Value = $u.D
    $ws.Range("J$r").Value = $u.J
    $ws.Range("K$r").Value = $u.K
    $ws.Range("L$r").Value = $u.L
    $ws.Range("M$r").Value = $u.M
    $ws.Range("P$r").Value = $u.P
}
